$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(59).Insert()
$ws.Range("A59").Value = 5
$ws.Range("B59").Value = "Macroferia Regional de Talca"
$ws.Range("C59").Value = "Maule"
$ws.Range("D59").Value = 44544
$ws.Range("E59").Value = 7
$ws.Range("F59").Value = 100112045
$ws.Range("G59").Value = "Zapallo"
$ws.Range("H59").Value = "Camote"
$ws.Range("I59").Value = "1a nueva(o)"
$ws.Range("J59").Value = 900
$ws.Range("K59").Value = 600
$ws.Range("L59").Value = 600
$ws.Range("M59").Value = 600
$ws.Range("N59").Value = "$/kilo (volumen en unidades)"
$ws.Range("O59").Value = "Región de O'Higgins"
$ws.Range("P59").Value = 600
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = "Hortaliza"
